$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Van-e akadály?" column (E) for the second week's rows (7-10):
# append the new obstacle noted by the team ("Tanár lehetetlen időkeretet adott")
$ws.Range("E7").Value = "Tanár lehetetlen időkeretet adott"
$ws.Range("E8").Value = "Nem volt elérhető a github repo, Tanár lehetetlen időkeretet adott"
$ws.Range("E9").Value = "Tanár lehetetlen időkeretet adott"
$ws.Range("E10").Value = "Tanár lehetetlen időkeretet adott"

# Reflect the final cursor/viewport position left after editing
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollColumn = 2
